$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 6

$ws.Range("D3").Value = 10.31
$ws.Range("F3").Value = 10.35

$ws.Range("C4").Value = 9.69
$ws.Range("E4").Value = 10.7
$ws.Range("F4").Value = 10.18

$ws.Range("D5").Value = 9.300000000000001
$ws.Range("F5").Value = 10.18
$ws.Range("I5").Value = 7.62

$ws.Range("C6").Value = 9.65
$ws.Range("D6").Value = 9.82
$ws.Range("E6").Value = 9.82
$ws.Range("G6").Value = 10.31
$ws.Range("H6").Value = 10.38

$ws.Range("F7").Value = 9.69

$ws.Range("B8").Value = 14
$ws.Range("F8").Value = 9.619999999999999

$ws.Range("E9").Value = 12.38
